$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.544602274894714
$ws.Range("B1").Value = 1.782876253128052
$ws.Range("C1").Value = 1.832174062728882
$ws.Range("D1").Value = 2.254343032836914
$ws.Range("E1").Value = 3.303630113601685
